$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header strings: volume/issue number and the reporting week date range.
$ws.Range("A8").Value = "Volume 30   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/4/2023  Through  12/10/2023"

# Updated weekly crime-complaint figures (new data pull for the week).
$ws.Range("L15").Value = 21.428571428571
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 25
$ws.Range("I16").Value = 124
$ws.Range("J16").Value = 129
$ws.Range("K16").Value = -3.875968992248
$ws.Range("L16").Value = -4.615384615384
$ws.Range("M16").Value = 39.325842696629
$ws.Range("N16").Value = -84.710234278668
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 500
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = 18.181818181818
$ws.Range("I17").Value = 138
$ws.Range("J17").Value = 127
$ws.Range("K17").Value = 8.661417322834
$ws.Range("L17").Value = 20
$ws.Range("M17").Value = 112.307692307692
$ws.Range("N17").Value = -20.231213872832
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = 21.428571428571
$ws.Range("I18").Value = 188
$ws.Range("J18").Value = 243
$ws.Range("K18").Value = -22.633744855967
$ws.Range("L18").Value = 8.045977011494
$ws.Range("M18").Value = 13.939393939393
$ws.Range("N18").Value = -77.960140679953
$ws.Range("C19").Value = 24
$ws.Range("D19").Value = 26
$ws.Range("E19").Value = -7.692307692307
$ws.Range("F19").Value = 95
$ws.Range("G19").Value = 99
$ws.Range("H19").Value = -4.040404040404
$ws.Range("I19").Value = 1222
$ws.Range("J19").Value = 1170
$ws.Range("K19").Value = 4.444444444444
$ws.Range("L19").Value = 53.517587939698
$ws.Range("M19").Value = 20.275590551181
$ws.Range("N19").Value = -66.793478260869
$ws.Range("C20").Value = 4
$ws.Range("E20").Value = 300
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 74
$ws.Range("J20").Value = 71
$ws.Range("K20").Value = 4.225352112676
$ws.Range("L20").Value = 12.121212121212
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = -91.405342624854
$ws.Range("C21").Value = 38
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = 26.666666666666
$ws.Range("F21").Value = 145
$ws.Range("G21").Value = 141
$ws.Range("H21").Value = 2.836879432624
$ws.Range("I21").Value = 1763
$ws.Range("J21").Value = 1761
$ws.Range("K21").Value = 0.113571834185
$ws.Range("L21").Value = 35.929067077872
$ws.Range("M21").Value = 27.939042089985
$ws.Range("N21").Value = -72.453125
$ws.Range("C22").Value = 5
$ws.Range("E22").Value = 150
$ws.Range("F22").Value = 11
$ws.Range("H22").Value = 83.333333333333
$ws.Range("I22").Value = 105
$ws.Range("J22").Value = 91
$ws.Range("K22").Value = 15.384615384615
$ws.Range("L22").Value = 34.615384615384
$ws.Range("M22").Value = 59.090909090909
$ws.Range("C24").Value = 73
$ws.Range("D24").Value = 77
$ws.Range("E24").Value = -5.194805194805
$ws.Range("F24").Value = 363
$ws.Range("G24").Value = 370
$ws.Range("H24").Value = -1.891891891891
$ws.Range("I24").Value = 3947
$ws.Range("J24").Value = 3806
$ws.Range("K24").Value = 3.704676826064
$ws.Range("L24").Value = 85.741176470588
$ws.Range("M24").Value = 143.943139678616
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = -45.454545454545
$ws.Range("F25").Value = 26
$ws.Range("G25").Value = 31
$ws.Range("H25").Value = -16.129032258064
$ws.Range("I25").Value = 389
$ws.Range("J25").Value = 354
$ws.Range("K25").Value = 9.887005649717
$ws.Range("L25").Value = 11.461318051575
$ws.Range("M25").Value = 58.775510204081
$ws.Range("L26").Value = 11.111111111111
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = -12.5
$ws.Range("I27").Value = 94
$ws.Range("J27").Value = 112
$ws.Range("K27").Value = -16.071428571428
$ws.Range("L27").Value = 17.5
$ws.Range("F30").Value = 3
$ws.Range("I30").Value = 18
$ws.Range("K30").Value = 20
$ws.Range("L30").Value = 28.571428571428
